$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row (row 14) for the "Serine" compound which has only parent
# isotopologues (0 13C, 0 15N) so the data correctly reflects a compound
# whose only labeled data is on the parent peak.
$ws.Range("A14").Value = "Serine"
$ws.Range("B14").Value = 105.093
$ws.Range("C14").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 445783904
$ws.Range("G14").Value = 452330528
$ws.Range("H14").Value = 460922144
$ws.Range("I14").Value = 460856768
$ws.Range("J14").Value = 240814112
$ws.Range("K14").Value = 238327808
$ws.Range("L14").Value = 236562352
$ws.Range("M14").Value = 238082256
$ws.Range("N14").Value = 0
$ws.Range("O14").Value = 0
$ws.Range("P14").Value = 0
$ws.Range("Q14").Value = 0

# Apply the same number-format style used in column B for other rows
$ws.Range("B14").NumberFormat = "0.0000"

$ws.Range("A14").Select()
